$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -12.2906
$ws.Range("C14").Value = -13.8451
$ws.Range("D15").Value = -7.8597
$ws.Range("C16").Value = -13.99339999999999
$ws.Range("C21").Value = -13.24770000000001
$ws.Range("D21").Value = -7.679499999999994
$ws.Range("D22").Value = -8.110200000000006
$ws.Range("C23").Value = -12.35140000000001
$ws.Range("D24").Value = -7.3935
$ws.Range("C25").Value = -11.8042
$ws.Range("C26").Value = -12.3348
$ws.Range("D27").Value = -7.980699999999998
$ws.Range("D28").Value = -7.869899999999994
$ws.Range("C29").Value = -11.02740000000001
$ws.Range("D36").Value = -6.8609
$ws.Range("D39").Value = -7.1172
$ws.Range("C40").Value = -11.8805
$ws.Range("D45").Value = -7.256800000000004
$ws.Range("D48").Value = -7.552799999999998
$ws.Range("D49").Value = -7.916300000000001
$ws.Range("D52").Value = -8.117200000000008
$ws.Range("C53").Value = -12.3767
$ws.Range("D53").Value = -8.696799999999994
$ws.Range("D54").Value = -7.789500000000003
$ws.Range("C57").Value = -14.04899999999999
$ws.Range("D57").Value = -7.921799999999998
$ws.Range("C59").Value = -12.553
$ws.Range("C65").Value = -12.6155
$ws.Range("C69").Value = -10.7644
$ws.Range("D70").Value = -6.920100000000001
$ws.Range("D71").Value = -7.327999999999996
$ws.Range("C79").Value = -11.6772
$ws.Range("C83").Value = -13.56019999999999
$ws.Range("D86").Value = -8.064199999999994
$ws.Range("D87").Value = -7.959200000000003
$ws.Range("D89").Value = -8.037799999999999
$ws.Range("C91").Value = -12.4858
$ws.Range("C93").Value = -10.15989999999999
$ws.Range("C100").Value = -12.1603
$ws.Range("D101").Value = -8.298200000000001
$ws.Range("C103").Value = -13.07379999999999
